# Updated scrape for diabetes one recipe.
# Fill in the "recipelist" sheet's header reorder/additions and the first
# data row (row 2) for "Chawli Bhaji".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("recipelist")

# --- Header row (row 1) ---
# Columns A-C unchanged. D/E/F are reordered (Ingredients, Preparation Time,
# Cooking Time now come before Food Category), G becomes Food Category, a new
# "Preparation Method" header is inserted at H, and I/J shift right by one
# (Nutrient values / Targetted morbid conditions) with a new "Recipe URL"
# header appended at K.
$ws.Range("D1").Value = "Ingredients"
$ws.Range("E1").Value = "Preparation Time"
$ws.Range("F1").Value = "Cooking Time"
$ws.Range("G1").Value = "Food Category"
$ws.Range("H1").Value = "Preparation Method"
$ws.Range("I1").Value = "Nutrient values"
$ws.Range("J1").Value = "Targetted morbid conditions"
$ws.Range("K1").Value = "Recipe URL"

# --- Data row (row 2) for the Chawli Bhaji recipe ---
$ws.Range("A2").Value = "Recipe# 6409`n04 Jan 19"
$ws.Range("B2").Value = "Chawli Bhaji"
$ws.Range("C2").Value = "Diabetic recipes recipes"
$ws.Range("D2").Value = "8 cups chopped chawli (cow pea / lobhia) leaves`n1/4 tsp turmeric powder (haldi)`nsalt to taste`n2 tsp oil`n1 tsp mustard seeds ( rai / sarson)`n8 to 10 curry leaves (kadi patta)`n4 whole dry kashmiri red chillies , broken into pieces`n2 tsp urad dal (split black lentils)`na pinch of asafoetida (hing)"
$ws.Range("E2").Value = "15 mins"
$ws.Range("F2").Value = "11 mins"
$ws.Range("G2").Value = "Vegetarian"
$ws.Range("H2").Value = "Method`nCombine the chawli leaves, turmeric powder, salt and 1¾ cups of water in a deep non-stick pan and cook on a medium flame for 5 to 7 minutes or till half of the water dries out. Keep aside to cool slightly.`nBlend in a mixer till smooth and keep aside.`nHeat the oil in a non-stick kadhai and add the mustard seeds.`nWhen the seeds crackle, add the curry leaves, red chillies, urad dal and asafoetida and sauté on a medium flame for a few seconds.`nAdd the chawli mixture and a little salt, mix well and cook on a medium flame for 2 to 3 minutes, while stirring occasionally.`nServe hot."
$ws.Range("I2").Value = "Accompaniments`nNutritious Lehsuni Methi Roti `nNutrient values (Abbrv) per serving`nEnergy 91 cal`nProtein 5.9 g`nCarbohydrates 9.5 g`nFiber 5.7 g`nFat 3.2 g`nCholesterol 0 mg`nSodium 313.6 mg`nClick here to view calories for Chawli Bhaji"
$ws.Range("J2").Value = "Diabetes"
$ws.Range("K2").Value = "https://www.tarladalal.com/chawli-bhaji-6409r"
